$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "p_adj" in H1
$ws.Range("H1").Value = "p_adj"

# Move the yellow highlight from column E (p) to column H (p_adj) for significant rows
$highlightRows = @(9, 13, 15, 16, 17, 18)
foreach ($r in $highlightRows) {
    $ws.Cells.Item($r, 5).ClearFormats()
}

# Fill in the p_adj values for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0.108933775858419
$ws.Range("H7").Value = 0.268345162545606
$ws.Range("H8").Value = 0.108933775858419
$ws.Range("H9").Value = 0.00539143282727617
$ws.Range("H10").Value = 0.511394742240951
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0.0284430238180771
$ws.Range("H14").Value = 0.108933775858419
$ws.Range("H15").Value = 0.00239816128374848
$ws.Range("H16").Value = 0.00239816128374848
$ws.Range("H17").Value = 0.0131589266691921
$ws.Range("H18").Value = 0.0189830314769883
$ws.Range("H19").Value = 1

# Re-apply yellow highlight fill to the new p_adj cells that correspond to significant p-values
foreach ($r in $highlightRows) {
    $ws.Cells.Item($r, 8).Interior.Color = 65535
}

# Restore the active selection to match the edited cell
$ws.Range("H9").Select()
